$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, copying the header style used by the other
# header cells (bold, centered, bordered) from G1 so it reuses the same
# style index rather than creating a duplicate one.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Fill in the new "Save" column values for the data rows.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
